$wb = $excel.ActiveWorkbook

# Sheet ALC, row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 689.7273
$ws.Range("I53").Value = 319.16666
$ws.Range("K53").Value = 319.16666
$ws.Range("M53").Value = 317.83334

# Sheet ALC, row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 683064.2
$ws.Range("I135").Value = 282.52777
$ws.Range("K135").Value = 2542.74993
$ws.Range("M135").Value = -7.749929999999949

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3481.8948
$ws.Range("I138").Value = 3859.6924
$ws.Range("J138").Value = 2663.3333
$ws.Range("K138").Value = 11579.0772
$ws.Range("L138").Value = 7989.999899999999
$ws.Range("M138").Value = -6439.0772
$ws.Range("N138").Value = -18269.9999

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4764.8154
$ws.Range("I32").Value = 4770.1763
$ws.Range("K32").Value = 4770.1763
$ws.Range("M32").Value = -4483.1763

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1815.9565
$ws.Range("I61").Value = 1743.8422
$ws.Range("J61").Value = 2158.5
$ws.Range("K61").Value = 1743.8422
$ws.Range("L61").Value = 2158.5
$ws.Range("M61").Value = -1531.8422
$ws.Range("N61").Value = -2582.5

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4174.421
$ws.Range("I74").Value = 714.7778
$ws.Range("K74").Value = 714.7778
$ws.Range("M74").Value = 159.2222

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4174.421
$ws.Range("I77").Value = 714.7778
$ws.Range("K77").Value = 3573.889
$ws.Range("M77").Value = 794.1110000000003

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1108
$ws.Range("I122").Value = 1184.5714
$ws.Range("J122").Value = 974
$ws.Range("K122").Value = 3553.7142
$ws.Range("L122").Value = 2922
$ws.Range("M122").Value = -1103.7142
$ws.Range("N122").Value = -7822

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 50398.906
$ws.Range("I132").Value = 2913.2856
$ws.Range("J132").Value = 145370.14
$ws.Range("K132").Value = 8739.856800000001
$ws.Range("L132").Value = 436110.42
$ws.Range("M132").Value = -6209.856800000001
$ws.Range("N132").Value = -441170.42

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1815.9565
$ws.Range("I136").Value = 1743.8422
$ws.Range("J136").Value = 2158.5
$ws.Range("K136").Value = 5231.5266
$ws.Range("L136").Value = 6475.5
$ws.Range("M136").Value = -2681.5266
$ws.Range("N136").Value = -11575.5

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1070
$ws.Range("I99").Value = 842
$ws.Range("J99").Value = 1640
$ws.Range("K99").Value = 842
$ws.Range("L99").Value = 1640
$ws.Range("M99").Value = 656
$ws.Range("N99").Value = -4636

# Sheet BSM, row 129
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H129").Value = 44110.89
$ws.Range("J129").Value = 44110.89
$ws.Range("L129").Value = 44110.89
$ws.Range("N129").Value = -54110.89

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1441.75
$ws.Range("I134").Value = 1482.8
$ws.Range("J134").Value = 1099.6666
$ws.Range("K134").Value = 4448.4
$ws.Range("L134").Value = 3298.9998
$ws.Range("M134").Value = -1913.4
$ws.Range("N134").Value = -8368.9998

# Sheet CUL, row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 7034.615
$ws.Range("I133").Value = 3330
$ws.Range("J133").Value = 9350
$ws.Range("K133").Value = 9990
$ws.Range("L133").Value = 28050
$ws.Range("M133").Value = -4930
$ws.Range("N133").Value = -38170

# Sheet GSM, row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2235.2666
$ws.Range("I126").Value = 2025.625
$ws.Range("J126").Value = 2474.8572
$ws.Range("K126").Value = 6076.875
$ws.Range("L126").Value = 7424.571599999999
$ws.Range("M126").Value = -3606.875
$ws.Range("N126").Value = -12364.5716

# Sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2657.1428
$ws.Range("I7").Value = 2533.3333
$ws.Range("K7").Value = 2533.3333
$ws.Range("M7").Value = -2421.3333

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2000.5714
$ws.Range("I40").Value = 2000.5714
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2000.5714
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1864.5714
$ws.Range("N40").ClearContents()

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1261.8
$ws.Range("I61").Value = 1239.5
$ws.Range("K61").Value = 1239.5
$ws.Range("M61").Value = -1037.5

# Sheet LTW, row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1972.6666
$ws.Range("I68").Value = 1791
$ws.Range("J68").Value = 2336
$ws.Range("K68").Value = 1791
$ws.Range("L68").Value = 2336
$ws.Range("M68").Value = -1042
$ws.Range("N68").Value = -3834

# Sheet LTW, row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1972.6666
$ws.Range("I71").Value = 1791
$ws.Range("J71").Value = 2336
$ws.Range("K71").Value = 8955
$ws.Range("L71").Value = 11680
$ws.Range("M71").Value = -5211
$ws.Range("N71").Value = -19168

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2529.7334
$ws.Range("I93").Value = 2557.3845
$ws.Range("J93").Value = 2350
$ws.Range("K93").Value = 2557.3845
$ws.Range("L93").Value = 2350
$ws.Range("M93").Value = -1309.3845
$ws.Range("N93").Value = -4846

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1261.8
$ws.Range("I113").Value = 1239.5
$ws.Range("K113").Value = 1239.5
$ws.Range("M113").Value = 930.5

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1623.6154
$ws.Range("I122").Value = 1001.6
$ws.Range("J122").Value = 2012.375
$ws.Range("K122").Value = 3004.8
$ws.Range("L122").Value = 6037.125
$ws.Range("M122").Value = -554.8000000000002
$ws.Range("N122").Value = -10937.125

# Sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2657.1428
$ws.Range("I126").Value = 2533.3333
$ws.Range("K126").Value = 7599.999899999999
$ws.Range("M126").Value = -5129.999899999999

# Sheet LTW, row 129
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H129").Value = 35953
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 35953
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 35953
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -45953

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2043.4517
$ws.Range("I136").Value = 967
$ws.Range("J136").Value = 4000.6365
$ws.Range("K136").Value = 2901
$ws.Range("L136").Value = 12001.9095
$ws.Range("M136").Value = -351
$ws.Range("N136").Value = -17101.9095

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 40498.168
$ws.Range("I62").Value = 39598
$ws.Range("J62").Value = 44999
$ws.Range("K62").Value = 39598
$ws.Range("L62").Value = 44999
$ws.Range("M62").Value = -38974
$ws.Range("N62").Value = -46247

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 40498.168
$ws.Range("I65").Value = 39598
$ws.Range("J65").Value = 44999
$ws.Range("K65").Value = 197990
$ws.Range("L65").Value = 224995
$ws.Range("M65").Value = -194870
$ws.Range("N65").Value = -231235

# Sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5027.077
$ws.Range("I81").Value = 1282.4
$ws.Range("J81").Value = 7367.5
$ws.Range("K81").Value = 2564.8
$ws.Range("L81").Value = 14735
$ws.Range("M81").Value = -1503.8
$ws.Range("N81").Value = -16857

# Sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5027.077
$ws.Range("I84").Value = 1282.4
$ws.Range("J84").Value = 7367.5
$ws.Range("K84").Value = 12824
$ws.Range("L84").Value = 73675
$ws.Range("M84").Value = -7520
$ws.Range("N84").Value = -84283

# Sheet WVR, row 129
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 40430
$ws.Range("J129").Value = 40430
$ws.Range("L129").Value = 40430
$ws.Range("N129").Value = -50430
